$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the numeric-looking / date-looking values as plain text (matches
# the source file, where every cell in the table -- including IDs,
# amounts and dates -- is stored as a shared string, not a real
# number/date). Columns B (description) and G (currency code) are
# already non-numeric text, so they don't need this.
$ws.Range("A4:A8").NumberFormat = "@"
$ws.Range("C4:F8").NumberFormat = "@"

# Row 4: was "Various paper supplies" (109566) record -> becomes the
# "Beverages and Catering" (290611 / USD) record that used to live at row 6
$ws.Range("A4").Value = "290611"
$ws.Range("B4").Value = "Beverages and Catering"
$ws.Range("C4").Value = "2017-08-09"
$ws.Range("D4").Value = "17159"
$ws.Range("E4").Value = "3431.8"
$ws.Range("F4").Value = "20590.8"
$ws.Range("G4").Value = "USD"

# Row 5: becomes the "Various paper supplies" (109566 / EUR) record
$ws.Range("A5").Value = "109566"
$ws.Range("B5").Value = "Various paper supplies"
$ws.Range("C5").Value = "2017-09-14"
$ws.Range("D5").Value = "136672"
$ws.Range("E5").Value = "27334.4"
$ws.Range("F5").Value = "164006"
$ws.Range("G5").Value = "EUR"

# Row 6: becomes the "Beverages and Catering" (463540 / EUR) record
$ws.Range("A6").Value = "463540"
$ws.Range("B6").Value = "Beverages and Catering"
$ws.Range("C6").Value = "2017-09-27"
$ws.Range("D6").Value = "136526"
$ws.Range("E6").Value = "27305.2"
$ws.Range("F6").Value = "163831"
$ws.Range("G6").Value = "EUR"

# Row 7: becomes a copy of the "Professional Services" (819413 / RON) record
$ws.Range("A7").Value = "819413"
$ws.Range("B7").Value = "Professional Services"
$ws.Range("C7").Value = "2017-11-24"
$ws.Range("D7").Value = "242624"
$ws.Range("E7").Value = "48524.8"
$ws.Range("F7").Value = "291149"
$ws.Range("G7").Value = "RON"

# Row 8: becomes a copy of the "Professional Services" (819413 / RON) record
$ws.Range("A8").Value = "819413"
$ws.Range("B8").Value = "Professional Services"
$ws.Range("C8").Value = "2017-11-24"
$ws.Range("D8").Value = "242624"
$ws.Range("E8").Value = "48524.8"
$ws.Range("F8").Value = "291149"
$ws.Range("G8").Value = "RON"

$wb.Save()
